$d = $word.ActiveDocument

# --- Edit 1: "covetous journey of information based " -> "desirable supply" + " of information based "
$rng = $d.Content
$found = $rng.Find.Execute("covetous journey of information based ", $true, $false, $false, $false, $false, $true, 1, $false, "desirable supply", 2)
if ($found) {
    $insertRng = $d.Range($rng.End, $rng.End)
    $insertRng.InsertAfter(" of information based ")
}

# --- Edit 2: "Austin – " -> "Feedback 1 – ..." + " for the reader. ..."
$rng = $d.Content
$found = $rng.Find.Execute("Austin – ", $true, $false, $false, $false, $false, $true, 1, $false, "Feedback 1 – It sounds great, it’s clear that the website offers great points", 2)
if ($found) {
    $insertRng = $d.Range($rng.End, $rng.End)
    $insertRng.InsertAfter(" for the reader. It has good details but some words could have been changed (the changes were already applied after feedback). ")
}

# --- Edit 3: "Natalie –" -> "Feedback 2 – " + "It seems to be well-written..."
$rng = $d.Content
$found = $rng.Find.Execute("Natalie –", $true, $false, $false, $false, $false, $true, 1, $false, "Feedback 2 – ", 2)
if ($found) {
    $insertRng = $d.Range($rng.End, $rng.End)
    $insertRng.InsertAfter("It seems to be well-written and brings the message across well, maybe the opening questions can be rephrased so they’re more direct. It’s important to eliminate redundant vocabulary and keep the message clear. It is overall still very clear but it can be a little more concise.")
}

# --- Edit 4: "Niza – " -> "Feedback 3 - " + "The elevator pitch..." + "enthusiastic;" + " it is easy to read..."
$rng = $d.Content
$found = $rng.Find.Execute("Niza – ", $true, $false, $false, $false, $false, $true, 1, $false, "Feedback 3 - ", 2)
if ($found) {
    $insertRng = $d.Range($rng.End, $rng.End)
    $insertRng.InsertAfter("The elevator pitch sounds very positive and ")

    $insertRng2 = $d.Range($insertRng.End, $insertRng.End)
    $insertRng2.InsertAfter("enthusiastic;")

    $insertRng3 = $d.Range($insertRng2.End, $insertRng2.End)
    $insertRng3.InsertAfter(" it is easy to read. Maybe it would help if you could add more details to the experiences that will be written about, such as main topics and quick bullet points to entice readers even more. A persuasive website goes a long way!")
}
